$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = 53377
$ws.Range("E2").Value = 1688
$ws.Range("F2").Value = 1688
$ws.Range("G2").Value = 1025
$ws.Range("H2").Value = 396
$ws.Range("I2").Value = 423
$ws.Range("J2").Value = -26
$ws.Range("K2").Value = 47492
$ws.Range("L2").Value = 27401
$ws.Range("M2").Value = 20091
$ws.Range("N2").Value = 18844
$ws.Range("O2").Value = 1247
$ws.Range("P2").Value = 1394
$ws.Range("Q2").Value = 1650
$ws.Range("R2").Value = -1752
$ws.Range("S2").Value = -521
$ws.Range("T2").Value = 2459
$ws.Range("U2").Value = -809
$ws.Range("V2").Value = 18442
$ws.Range("W2").Value = 3.16
$ws.Range("X2").Value = 0.74
$ws.Range("Y2").Value = 2.28
$ws.Range("Z2").Value = 0.84
$ws.Range("AA2").Value = 136.38
$ws.Range("AB2").Value = 1289.76
$ws.Range("AC2").Value = 1517
$ws.Range("AD2").Value = 31.87
$ws.Range("AE2").Value = 67609
$ws.Range("AF2").Value = 0.72
$ws.Range("AG2").Value = 500
$ws.Range("AH2").Value = 1.03
$ws.Range("AI2").Value = 33.29
$ws.Range("AJ2").Value = 25103951
$ws.Range("D3").Value = 48565
$ws.Range("E3").Value = 2805
$ws.Range("F3").Value = 2805
$ws.Range("G3").Value = -1377
$ws.Range("H3").Value = -1451
$ws.Range("I3").Value = -1443
$ws.Range("J3").Value = -8
$ws.Range("K3").Value = 49861
$ws.Range("L3").Value = 30217
$ws.Range("M3").Value = 19644
$ws.Range("N3").Value = 18320
$ws.Range("O3").Value = 1324
$ws.Range("P3").Value = 1394
$ws.Range("Q3").Value = 3658
$ws.Range("R3").Value = -3766
$ws.Range("S3").Value = 411
$ws.Range("T3").Value = 2017
$ws.Range("U3").Value = 1640
$ws.Range("V3").Value = 18843
$ws.Range("W3").Value = 5.78
$ws.Range("X3").Value = -2.99
$ws.Range("Y3").Value = -7.77
$ws.Range("Z3").Value = -2.98
$ws.Range("AA3").Value = 153.82
$ws.Range("AB3").Value = 1182.81
$ws.Range("AC3").Value = -5177
$ws.Range("AD3").Value = -12.13
$ws.Range("AE3").Value = 65694
$ws.Range("AF3").Value = 0.96
$ws.Range("AG3").Value = 500
$ws.Range("AH3").Value = 0.8
$ws.Range("AI3").Value = -9.76
$ws.Range("AJ3").Value = 25119218
$ws.Range("D4").Value = 45622
$ws.Range("E4").Value = 2767
$ws.Range("F4").Value = 2767
$ws.Range("G4").Value = 2389
$ws.Range("H4").Value = 1760
$ws.Range("I4").Value = 1702
$ws.Range("J4").Value = 57
$ws.Range("K4").Value = 52491
$ws.Range("L4").Value = 31239
$ws.Range("M4").Value = 21253
$ws.Range("N4").Value = 19810
$ws.Range("O4").Value = 1443
$ws.Range("P4").Value = 1396
$ws.Range("Q4").Value = 3396
$ws.Range("R4").Value = -2618
$ws.Range("S4").Value = -185
$ws.Range("T4").Value = 2431
$ws.Range("U4").Value = 965
$ws.Range("V4").Value = 19192
$ws.Range("W4").Value = 6.07
$ws.Range("X4").Value = 3.86
$ws.Range("Y4").Value = 8.93
$ws.Range("Z4").Value = 3.44
$ws.Range("AA4").Value = 146.99
$ws.Range("AB4").Value = 1282.74
$ws.Range("AC4").Value = 6102
$ws.Range("AD4").Value = 12.13
$ws.Range("AE4").Value = 70954
$ws.Range("AF4").Value = 1.04
$ws.Range("AG4").Value = 1100
$ws.Range("AH4").Value = 1.49
$ws.Range("AI4").Value = 18.12
$ws.Range("AJ4").Value = 25151485
$ws.Range("D5").Value = 46070
$ws.Range("E5").Value = 1980
$ws.Range("F5").Value = 1980
$ws.Range("G5").Value = 1793
$ws.Range("H5").Value = 1227
$ws.Range("I5").Value = 1351
$ws.Range("J5").Value = -123
$ws.Range("K5").Value = 56619
$ws.Range("L5").Value = 33395
$ws.Range("M5").Value = 23224
$ws.Range("N5").Value = 21833
$ws.Range("O5").Value = 1391
$ws.Range("P5").Value = 1421
$ws.Range("Q5").Value = 2766
$ws.Range("R5").Value = -4987
$ws.Range("S5").Value = 1576
$ws.Range("T5").Value = 4548
$ws.Range("U5").Value = -1782
$ws.Range("V5").Value = 21022
$ws.Range("W5").Value = 4.3
$ws.Range("X5").Value = 2.66
$ws.Range("Y5").Value = 6.49
$ws.Range("Z5").Value = 2.25
$ws.Range("AA5").Value = 143.79
$ws.Range("AB5").Value = 1340.69
$ws.Range("AC5").Value = 4830
$ws.Range("AD5").Value = 18.43
$ws.Range("AE5").Value = 76824
$ws.Range("AF5").Value = 1.16
$ws.Range("AG5").Value = 1100
$ws.Range("AH5").Value = 1.24
$ws.Range("AI5").Value = 23.25
$ws.Range("AJ5").Value = 25651728
$ws.Range("D6").Value = 47526
$ws.Range("E6").Value = 1449
$ws.Range("F6").Value = 1449
$ws.Range("G6").Value = 903
$ws.Range("H6").Value = 429
$ws.Range("I6").Value = 650
$ws.Range("K6").Value = 55729
$ws.Range("L6").Value = 33643
$ws.Range("M6").Value = 22085
$ws.Range("N6").Value = 20965
$ws.Range("P6").Value = 1487
$ws.Range("Q6").Value = 3412
$ws.Range("R6").Value = -2959
$ws.Range("S6").Value = -361
$ws.Range("T6").Value = 3832
$ws.Range("U6").Value = -420
$ws.Range("V6").Value = 22672
$ws.Range("W6").Value = 3.05
$ws.Range("X6").Value = 0.9
$ws.Range("Y6").Value = 3.04
$ws.Range("Z6").Value = 0.76
$ws.Range("AA6").Value = 152.33
$ws.Range("AB6").Value = 1330.71
$ws.Range("AC6").Value = 2250
$ws.Range("AD6").Value = 25.47
$ws.Range("AE6").Value = 70480
$ws.Range("AF6").Value = 0.8100000000000001
$ws.Range("AG6").Value = 900
$ws.Range("AH6").Value = 1.57
$ws.Range("AI6").Value = 41.38
$ws.Range("AJ6").Value = 26978840
$ws.Range("D7").Value = 45046
$ws.Range("E7").Value = 2222
$ws.Range("G7").Value = 1889
$ws.Range("H7").Value = 1202
$ws.Range("I7").Value = 1275
$ws.Range("K7").Value = 55520
$ws.Range("L7").Value = 33226
$ws.Range("M7").Value = 22295
$ws.Range("N7").Value = 21314
$ws.Range("P7").Value = 1488
$ws.Range("Q7").Value = 3830
$ws.Range("R7").Value = -812
$ws.Range("S7").Value = -1451
$ws.Range("T7").Value = 1704
$ws.Range("U7").Value = 2231
$ws.Range("W7").Value = 4.93
$ws.Range("X7").Value = 2.67
$ws.Range("Y7").Value = 6.03
$ws.Range("Z7").Value = 2.16
$ws.Range("AA7").Value = 149.03
$ws.Range("AC7").Value = 4286
$ws.Range("AD7").Value = 9.859999999999999
$ws.Range("AE7").Value = 71652
$ws.Range("AF7").Value = 0.59
$ws.Range("AG7").Value = 952
$ws.Range("AH7").Value = 2.25
$ws.Range("AI7").Value = 20.14
$ws.Range("D8").Value = 47779
$ws.Range("E8").Value = 2475
$ws.Range("G8").Value = 2317
$ws.Range("H8").Value = 1758
$ws.Range("I8").Value = 1789
$ws.Range("K8").Value = 56516
$ws.Range("L8").Value = 33005
$ws.Range("M8").Value = 23511
$ws.Range("N8").Value = 22236
$ws.Range("P8").Value = 1488
$ws.Range("Q8").Value = 4162
$ws.Range("R8").Value = -2044
$ws.Range("S8").Value = -1827
$ws.Range("T8").Value = 2159
$ws.Range("U8").Value = 1372
$ws.Range("W8").Value = 5.18
$ws.Range("X8").Value = 3.68
$ws.Range("Y8").Value = 8.210000000000001
$ws.Range("Z8").Value = 3.14
$ws.Range("AA8").Value = 140.38
$ws.Range("AC8").Value = 6014
$ws.Range("AD8").Value = 7.03
$ws.Range("AE8").Value = 74750
$ws.Range("AF8").Value = 0.57
$ws.Range("AG8").Value = 1033
$ws.Range("AH8").Value = 2.45
$ws.Range("AI8").Value = 15.58
$ws.Range("D9").Value = 49544
$ws.Range("E9").Value = 2720
$ws.Range("G9").Value = 2258
$ws.Range("H9").Value = 1728
$ws.Range("I9").Value = 1765
$ws.Range("K9").Value = 57792
$ws.Range("L9").Value = 33206
$ws.Range("M9").Value = 24586
$ws.Range("N9").Value = 23288
$ws.Range("P9").Value = 1488
$ws.Range("Q9").Value = 4321
$ws.Range("R9").Value = -1999
$ws.Range("S9").Value = -1298
$ws.Range("T9").Value = 2148
$ws.Range("U9").Value = 1291
$ws.Range("W9").Value = 5.49
$ws.Range("X9").Value = 3.49
$ws.Range("Y9").Value = 7.75
$ws.Range("Z9").Value = 3.02
$ws.Range("AA9").Value = 135.06
$ws.Range("AC9").Value = 5933
$ws.Range("AD9").Value = 7.12
$ws.Range("AE9").Value = 78288
$ws.Range("AF9").Value = 0.54
$ws.Range("AG9").Value = 1067
$ws.Range("AH9").Value = 16.3
